# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F (dSF) after repulling data
$ws.Range("F8").Value = 1
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 1
$ws.Range("F36").Value = 2
$ws.Range("F38").Value = 1
$ws.Range("F43").Value = 1
$ws.Range("F47").Value = 4
$ws.Range("F50").Value = -3
$ws.Range("F53").Value = -4
$ws.Range("F61").Value = -2
$ws.Range("F63").Value = -2
$ws.Range("F64").Value = -2
